$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.552.01"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "2.616.00"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'534.07"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'142.68"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "2.619.43"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").Value = "3.075.81"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "58.516.51"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "'20.76"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "2.601.40"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'334.47"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "'6.21"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'66.68"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "'0.420"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -2.77%  "
$ws.Range("D28").Value = "'7.09"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Value = "0.0₃0733"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.01"
$ws.Range("E31").Value = "  +3.26%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.63"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "'154.74"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").Value = "'18.94"
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").Value = "'0.836"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'0.815"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.42"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").Value = "'286.38"
$ws.Range("E41").Value = "  +4.17%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'0.596"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'0.0944"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "'19.04"
$ws.Range("E46").Value = "  +2.71%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").Value = "'0.0225"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "1.938.43"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'17.84"
$ws.Range("E51").Value = "  -3.04%  "
